$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 418, shifting existing rows 418-460 down to 419-461.
$ws.Rows.Item(418).Insert()

# Populate the new row 418 with a fresh weekly record.
# Columns A,B,C,E,F,G,H,I,O,R are identical to the row that used to be
# at 418 (now at 419), so Excel's default "format/copy from row above"
# behavior combined with explicit values below reproduces them.
$ws.Cells.Item(418, 1).Value = 10
$ws.Cells.Item(418, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(418, 3).Value = "La Araucanía"
$ws.Cells.Item(418, 4).Value = 44826
$ws.Cells.Item(418, 5).Value = 9
$ws.Cells.Item(418, 6).Value = 100114014
$ws.Cells.Item(418, 7).Value = "Betarraga"
$ws.Cells.Item(418, 8).Value = "Sin especificar"
$ws.Cells.Item(418, 9).Value = "Primera"
$ws.Cells.Item(418, 10).Value = 20
$ws.Cells.Item(418, 11).Value = 10000
$ws.Cells.Item(418, 12).Value = 10000
$ws.Cells.Item(418, 13).Value = 10000
$ws.Cells.Item(418, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(418, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(418, 16).Value = 400
$ws.Cells.Item(418, 17).Value = 25
$ws.Cells.Item(418, 18).Value = "Hortaliza"
